# Davuluri_LabExam03Grading.xlsx - grading update
# - Question 12 (row 20): score dropped from 10 to 9, with a new grading
#   comment explaining the deduction.
# - Question 18 / CustomerMappingTest Class (row 34): score raised from 0
#   to 3, with the old terse comment replaced by a longer explanation
#   (wrapped, taller row).
# Dependent subtotal/total formulas recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 20: addProduct() method grading ---
$ws.Range("E20").Value = 9
$ws.Range("F20").Value = "(-1)For incorrect condition for checking customer exists or not."

# --- Row 34: CustomerMappingTest Class grading ---
$ws.Range("E34").Value = 3

$rightQuote = [char]0x2019
$f34Comment = "(-4)I have changed your addProduct() code and run the test cases then 4 test cases failed but I didn" + $rightQuote + "t deducted any points for remaining test cases`n"
$ws.Range("F34").Value = $f34Comment
$ws.Range("F34").WrapText = $true
$ws.Rows.Item(34).RowHeight = 30

# --- Move the on-screen selection/scroll position to the area just edited ---
$ws.Activate()
$ws.Range("F34").Select()
$excel.ActiveWindow.ScrollRow = 20
